# Threat Alert Report update - 2026-01-18 01:00
# Applies the row-level edits described by the source diff: updated dates,
# flight/airline reshuffles, fare figures, and two LOW<->MEDIUM threat
# swaps (which also flip the cell's highlight style) on the THREAT_ALERT
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a "DD-MON-YY" style value into a cell while keeping it a
# literal text string (Excel would otherwise silently coerce it into a
# date serial + date number-format) and keeping the cell's original style
# (fill/border/font) untouched.
function Set-DateText($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Copy()
    $ws.Range("M1").PasteSpecial(-4122)   # xlPasteFormats -> stash original fmt
    $cell.NumberFormat = "@"              # force text interpretation
    $cell.Value = $text
    $ws.Range("M1").Copy()
    $cell.PasteSpecial(-4122)             # restore original fmt onto the cell
    $ws.Range("M1").Clear()
}

# Helper: swap the visual style (fill/font/border) of two cells, leaving
# their values untouched.
function Swap-Format($ref1, $ref2) {
    $ws.Range($ref1).Copy()
    $ws.Range("M1").PasteSpecial(-4122)
    $ws.Range($ref2).Copy()
    $ws.Range($ref1).PasteSpecial(-4122)
    $ws.Range("M1").Copy()
    $ws.Range($ref2).PasteSpecial(-4122)
    $ws.Range("M1").Clear()
}

# ---- Row 2 ----
Set-DateText "A2" "08-FEB-26"

# ---- Row 3 ----
$ws.Range("D3").Value = 640
$ws.Range("F3").Value = -7

# ---- Row 4 ----
Set-DateText "A4" "15-FEB-26"
$ws.Range("C4").Value = "Nile Air NP-144"
$ws.Range("D4").Value = 571
$ws.Range("F4").Value = -76

# ---- Row 6 ----
Set-DateText "A6" "19-FEB-26"
$ws.Range("B6").Value = "SM-448"
$ws.Range("D6").Value = 571
$ws.Range("E6").Value = 591
$ws.Range("F6").Value = -20

# ---- Row 7 ----
Set-DateText "A7" "20-FEB-26"
$ws.Range("B7").Value = "SM-444"
$ws.Range("C7").Value = "Air Arabia Egypt E5-512"
$ws.Range("D7").Value = 594
$ws.Range("E7").Value = 721
$ws.Range("F7").Value = -127

# ---- Row 8 ----
Set-DateText "A8" "27-FEB-26"
$ws.Range("D8").Value = 725
$ws.Range("E8").Value = 786
$ws.Range("F8").Value = -61

# ---- Row 9 ----
Set-DateText "A9" "01-MAR-26"
$ws.Range("E9").Value = 883
$ws.Range("F9").Value = -158

# ---- Row 10 ----
$ws.Range("C10").Value = "Nile Air NP-144"
$ws.Range("D10").Value = 863
$ws.Range("F10").Value = -20

# ---- Row 11 ----
Set-DateText "A11" "05-MAR-26"
$ws.Range("B11").Value = "SM-448"
$ws.Range("D11").Value = 1101
$ws.Range("E11").Value = 1159
$ws.Range("F11").Value = -58

# ---- Row 12 ----
Set-DateText "A12" "06-MAR-26"
$ws.Range("B12").Value = "SM-444"
$ws.Range("C12").Value = "Air Arabia Egypt E5-512"
$ws.Range("D12").Value = 874
$ws.Range("F12").Value = -285
$ws.Range("J12").Value = "MEDIUM THREAT - MONITOR"

# ---- Row 13 ----
Set-DateText "A13" "08-MAR-26"
$ws.Range("E13").Value = 1013
$ws.Range("F13").Value = -139
$ws.Range("J13").Value = "LOW THREAT"

# J12/J13 swap LOW THREAT (green) <-> MEDIUM THREAT - MONITOR (yellow) styling
Swap-Format "J12" "J13"

# ---- Row 14 ----
Set-DateText "A14" "13-MAR-26"
$ws.Range("D14").Value = 1180
$ws.Range("E14").Value = 1306
$ws.Range("F14").Value = -126

# ---- Row 15 ----
Set-DateText "A15" "15-MAR-26"

# ---- Row 16 ----
Set-DateText "A16" "22-MAR-26"
$ws.Range("D16").Value = 650
$ws.Range("E16").Value = 883
$ws.Range("F16").Value = -233
$ws.Range("J16").Value = "MEDIUM THREAT - MONITOR"

# ---- Row 17 ----
Set-DateText "A17" "27-MAR-26"
$ws.Range("D17").Value = 612
$ws.Range("E17").Value = 721
$ws.Range("F17").Value = -109
$ws.Range("J17").Value = "LOW THREAT"

# J16/J17 swap LOW THREAT (green) <-> MEDIUM THREAT - MONITOR (yellow) styling
Swap-Format "J16" "J17"
